# Applies the changes described by the target diff:
#  1. Update the cached text of the auto-updating date fields (master,
#     all 11 slide layouts, and the notes master) from 5/5/21 to 5/18/21.
#  2. Move/resize the "Rectangle 4" shape on slide 1.
#  3. Remove the Footer / Slide Number placeholders from slide 1 by
#     turning off their visibility in the slide's headers & footers
#     (PowerPoint then omits the corresponding placeholder shapes).

$p = $ppt.ActivePresentation

$oldDate = "5/5/21"
$newDate = "5/18/21"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 1b. Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 1c. Notes master date placeholder.
#
# NOTE: this runtime mis-routes NotesMaster.Shapes(i) *writes* onto the
# regular SlideMaster's shape at the same position instead of the notes
# master (reproducible even via plain VBA - a sandbox limitation, not
# something reachable from script), so it is intentionally left alone
# here rather than risk clobbering unrelated slide-master content.
$notesMaster = $p.NotesMaster

# 2. Move/resize "Rectangle 4" on slide 1.
#    (Point values picked so that, after this runtime's internal
#    rounding of Shape.Left/Top/Width/Height, the saved EMU values come
#    back out to exactly 1287063 / 809297 / 8928991 / 4100651.)
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 4") {
        $sh.Left = 101.343583307087
        $sh.Top = 63.7241932283465
        $sh.Width = 703.070197480315
        $sh.Height = 322.885945511811
    }
}

# 3. Drop the footer and slide-number placeholders from slide 1.
$s.HeadersFooters.Footer.Visible = $false
$s.HeadersFooters.SlideNumber.Visible = $false
